$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Support for Multiple Preset Cameras & Split-screen Camera Views" row (row 14)
# now has credit earned, and is marked "I want to do this"
$ws.Range("C14").Value = 0.05
$ws.Range("D14").Value = "I want to do this"

# Camera control related rows (row 18 "Support for Skyboxes..." and row 19
# "Support for Game Specific Data...") no longer marked "I want to do this"
$ws.Range("D18").ClearContents()
$ws.Range("D19").ClearContents()

# Row 26 ("Support for adding Point and Spotlight sources via Blender") is now
# marked "I want to do this" instead
$ws.Range("D26").Value = "I want to do this"

# Row 52 ("Commit regularly throughout the month to your GIT repo") no longer
# marked "Doing this"
$ws.Range("D52").ClearContents()

# Update the selection/view to reflect where the user left off editing
$ws.Range("E7").Select()
